$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Fill in the new use case row (row 26): ID, Alias, Descripcion, Estado, Esfuerzo, Incremento, Prioridad
# (order of assignment matters for shared-string table ordering)
$ws.Range("B26").Value = "CU - 22"
$ws.Range("D26").Value = "Inciar sesion"
$ws.Range("C26").Value = "En este caso de uso el usuario del sistema ingresa a con algun tipo de autenticacion para poder ver y modificar sus datos."
$ws.Range("E26").Value = "vacio"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 1

# Row grows to match the wrapped text of the new description (same as other rows)
$ws.Rows.Item(26).RowHeight = 30

# Update the selection/view state to match the authored commit
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C23").Select()
